$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Última actualización: 06:17:28"
$ws1.Range("A3").Value = "Total filas: 25"
$ws1.Cells.Item(9,1).Value = "06:17:28"
$ws1.Cells.Item(9,2).Value = "06:17"
$ws1.Cells.Item(9,3).Value = "215A_EL PATO"
$ws1.Cells.Item(9,4).Value = 0
$ws1.Cells.Item(9,5).Value = "LP1912"
$ws1.Cells.Item(10,1).Value = "05:57:13"
$ws1.Cells.Item(10,2).Value = "06:30"
$ws1.Cells.Item(10,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(10,4).Value = 33
$ws1.Cells.Item(10,5).Value = "LP1912"
$ws1.Cells.Item(11,1).Value = "06:17:28"
$ws1.Cells.Item(11,2).Value = "06:32"
$ws1.Cells.Item(11,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(11,4).Value = 15
$ws1.Cells.Item(11,5).Value = "LP1912"
$ws1.Cells.Item(12,1).Value = "05:57:13"
$ws1.Cells.Item(12,2).Value = "06:34"
$ws1.Cells.Item(12,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(12,4).Value = 37
$ws1.Cells.Item(12,5).Value = "LP1912"
$ws1.Cells.Item(13,1).Value = "05:57:13"
$ws1.Cells.Item(13,2).Value = "06:39"
$ws1.Cells.Item(13,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(13,4).Value = 42
$ws1.Cells.Item(13,5).Value = "LP1912"
$ws1.Cells.Item(14,1).Value = "05:57:13"
$ws1.Cells.Item(14,2).Value = "06:41"
$ws1.Cells.Item(14,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(14,4).Value = 44
$ws1.Cells.Item(14,5).Value = "LP1912"
$ws1.Cells.Item(15,1).Value = "05:57:13"
$ws1.Cells.Item(15,2).Value = "06:57"
$ws1.Cells.Item(15,3).Value = "215A_EL PATO"
$ws1.Cells.Item(15,4).Value = 60
$ws1.Cells.Item(15,5).Value = "LP1912"
$ws1.Cells.Item(16,1).Value = "05:57:13"
$ws1.Cells.Item(16,2).Value = "06:59"
$ws1.Cells.Item(16,3).Value = "225_GOMEZ"
$ws1.Cells.Item(16,4).Value = 62
$ws1.Cells.Item(16,5).Value = "LP1912"
$ws1.Cells.Item(17,1).Value = "06:17:28"
$ws1.Cells.Item(17,2).Value = "07:15"
$ws1.Cells.Item(17,3).Value = "215C_EL PATO"
$ws1.Cells.Item(17,4).Value = 58
$ws1.Cells.Item(17,5).Value = "LP1912"
$ws1.Cells.Item(18,1).Value = "05:57:13"
$ws1.Cells.Item(18,2).Value = "07:16"
$ws1.Cells.Item(18,3).Value = "215C_EL PATO"
$ws1.Cells.Item(18,4).Value = 79
$ws1.Cells.Item(18,5).Value = "LP1912"
$ws1.Cells.Item(19,1).Value = "05:57:13"
$ws1.Cells.Item(19,2).Value = "07:19"
$ws1.Cells.Item(19,3).Value = "14_ABASTO"
$ws1.Cells.Item(19,4).Value = 82
$ws1.Cells.Item(19,5).Value = "LP1912"
$ws1.Cells.Item(20,1).Value = "05:57:13"
$ws1.Cells.Item(20,2).Value = "07:21"
$ws1.Cells.Item(20,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(20,4).Value = 84
$ws1.Cells.Item(20,5).Value = "LP1912"
$ws1.Cells.Item(21,1).Value = "06:17:28"
$ws1.Cells.Item(21,2).Value = "07:21"
$ws1.Cells.Item(21,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(21,4).Value = 64
$ws1.Cells.Item(21,5).Value = "LP1912"
$ws1.Cells.Item(22,1).Value = "05:57:13"
$ws1.Cells.Item(22,2).Value = "07:29"
$ws1.Cells.Item(22,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(22,4).Value = 92
$ws1.Cells.Item(22,5).Value = "LP1912"
$ws1.Cells.Item(23,1).Value = "05:57:13"
$ws1.Cells.Item(23,2).Value = "07:35"
$ws1.Cells.Item(23,3).Value = "10_OLMOS"
$ws1.Cells.Item(23,4).Value = 98
$ws1.Cells.Item(23,5).Value = "LP1912"
$ws1.Cells.Item(24,1).Value = "06:17:28"
$ws1.Cells.Item(24,2).Value = "07:36"
$ws1.Cells.Item(24,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(24,4).Value = 79
$ws1.Cells.Item(24,5).Value = "LP1912"
$ws1.Cells.Item(25,1).Value = "05:57:13"
$ws1.Cells.Item(25,2).Value = "07:37"
$ws1.Cells.Item(25,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(25,4).Value = 100
$ws1.Cells.Item(25,5).Value = "LP1912"
$ws1.Cells.Item(26,1).Value = "05:57:13"
$ws1.Cells.Item(26,2).Value = "07:55"
$ws1.Cells.Item(26,3).Value = "14_ABASTO"
$ws1.Cells.Item(26,4).Value = 118
$ws1.Cells.Item(26,5).Value = "LP1912"
$ws1.Cells.Item(27,1).Value = "06:17:28"
$ws1.Cells.Item(27,2).Value = "08:00"
$ws1.Cells.Item(27,3).Value = "17_ROMERO"
$ws1.Cells.Item(27,4).Value = 103
$ws1.Cells.Item(27,5).Value = "LP1912"
$ws1.Cells.Item(28,1).Value = "06:17:28"
$ws1.Cells.Item(28,2).Value = "08:01"
$ws1.Cells.Item(28,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(28,4).Value = 104
$ws1.Cells.Item(28,5).Value = "LP1912"
$ws1.Cells.Item(29,1).Value = "06:17:28"
$ws1.Cells.Item(29,2).Value = "08:11"
$ws1.Cells.Item(29,3).Value = "10_OLMOS"
$ws1.Cells.Item(29,4).Value = 114
$ws1.Cells.Item(29,5).Value = "LP1912"
$ws1.Cells.Item(30,1).Value = "06:17:28"
$ws1.Cells.Item(30,2).Value = "08:13"
$ws1.Cells.Item(30,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(30,4).Value = 116
$ws1.Cells.Item(30,5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 06:17:28"
$ws2.Range("A3").Value = "Total filas: 5"
$ws2.Cells.Item(7,1).Value = "06:17:28"
$ws2.Cells.Item(7,2).Value = "06:17"
$ws2.Cells.Item(7,3).Value = "215A_EL PATO"
$ws2.Cells.Item(7,4).Value = 0
$ws2.Cells.Item(7,5).Value = "LP1912"
$ws2.Cells.Item(8,1).Value = "06:17:28"
$ws2.Cells.Item(8,2).Value = "06:57"
$ws2.Cells.Item(8,3).Value = "215A_EL PATO"
$ws2.Cells.Item(8,4).Value = 60
$ws2.Cells.Item(8,5).Value = "LP1912"
$ws2.Cells.Item(9,1).Value = "06:17:28"
$ws2.Cells.Item(9,2).Value = "07:15"
$ws2.Cells.Item(9,3).Value = "215C_EL PATO"
$ws2.Cells.Item(9,4).Value = 58
$ws2.Cells.Item(9,5).Value = "LP1912"
$ws2.Cells.Item(10,1).Value = "05:57:13"
$ws2.Cells.Item(10,2).Value = "07:16"
$ws2.Cells.Item(10,3).Value = "215C_EL PATO"
$ws2.Cells.Item(10,4).Value = 79
$ws2.Cells.Item(10,5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 06:17:28"
